# Normalize stray-quote / typo values in the "purpose" (C), "age" (G) and
# "job" (H) columns of the dataset, and rename a couple of categorical
# labels ("Old" -> "Senior", "the" -> "other").
#
# Note: Excel's Range.Value setter (and the Replacement side of
# Range.Replace) applies the classic "leading apostrophe = force-text"
# entry convention, silently swallowing a single leading ' from whatever
# is written into a cell. Values that must legitimately *start* with a
# literal apostrophe therefore need that apostrophe doubled in the
# PowerShell literal passed to Replace (''text -> stored as 'text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# --- Column C: purpose -----------------------------------------------
$colC = $ws.Range("C1:C$lastRow")

# Strip the (erroneous) surrounding quote marks that were stored
# literally in the cell text.
$colC.Replace("'new car'", "new car", 1, 1, $true)
$colC.Replace("'used car'", "used car", 1, 1, $true)
$colC.Replace("'domestic appliance'", "domestic appliance", 1, 1, $true)

# Fix entries that were missing their leading quote (typo in source data)
# so they end up with the same clean, quote-free text as the rest.
$colC.Replace("new car'", "new car", 1, 1, $true)
$colC.Replace("use car'", "used car", 1, 1, $true)

# Fix a couple of truncated/garbled "other" values.
$colC.Replace("the", "other", 1, 1, $true)

# --- Column G: age -----------------------------------------------------
$colG = $ws.Range("G1:G$lastRow")
$colG.Replace("Old", "Senior", 1, 1, $true)

# --- Column H: job -------------------------------------------------------
# A couple of rows are missing the quote marks that every other
# "unskilled resident" entry carries; add them back for consistency.
# (Leading apostrophe doubled so it survives the force-text stripping.)
$colH = $ws.Range("H1:H$lastRow")
$colH.Replace("unskilled resident", "''unskilled resident'", 1, 1, $true)
